$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "18 x 45" + [char]11 + "  4    5" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "8|    |"
$t.Cell(1,2).Range.Text = "78 x 37" + [char]11 + "  3    7" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "8|    |"
$t.Cell(1,3).Range.Text = "28 x 54" + [char]11 + "  5    4" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "8|    |"

$t.Cell(2,1).Range.Text = "34 x 13" + [char]11 + "  1    3" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "4|    |"
$t.Cell(2,2).Range.Text = "81 x 41" + [char]11 + "  4    1" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "1|    |"
$t.Cell(2,3).Range.Text = "35 x 88" + [char]11 + "  8    8" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "5|    |"

$t.Cell(3,1).Range.Text = "55 x 49" + [char]11 + "  4    9" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "5|    |"
$t.Cell(3,2).Range.Text = "97 x 40" + [char]11 + "  4    0" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
$t.Cell(3,3).Range.Text = "95 x 75" + [char]11 + "  7    5" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "5|    |"

$t.Cell(4,1).Range.Text = "24 x 43" + [char]11 + "  4    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "4|    |"
$t.Cell(4,2).Range.Text = "97 x 20" + [char]11 + "  2    0" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
$t.Cell(4,3).Range.Text = "79 x 30" + [char]11 + "  3    0" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "9|    |"

$t.Cell(5,1).Range.Text = "16 x 61" + [char]11 + "  6    1" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "6|    |"
$t.Cell(5,2).Range.Text = "17 x 36" + [char]11 + "  3    6" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "7|    |"
$t.Cell(5,3).Range.Text = "43 x 27" + [char]11 + "  2    7" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "3|    |"
